# Applies hybrid bold + color ("2C3E50") highlighting to quantitative
# metrics (percentages, dollar amounts, etc.) inside specific bullet /
# achievement paragraphs of the resume, matching the target diff. Each
# targeted number (e.g. "23%", "$4.7M", "±4.2%") is split out into its
# own run with Bold + Color formatting, while the surrounding text stays
# in plain runs.

$d = $word.ActiveDocument

# Word's Font.Color takes a BGR-packed long, so the byte order is
# reversed relative to the OOXML w:color hex value "2C3E50".
$highlightColor = 5258796

function Find-ParagraphIndex($UniqueSubstring) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text -like "*$UniqueSubstring*") {
            return $i
        }
    }
    throw "No paragraph found containing '$UniqueSubstring'"
}

function Highlight-InParagraph($UniqueSubstring, $Targets) {
    $idx = Find-ParagraphIndex $UniqueSubstring
    $p = $d.Paragraphs.Item($idx)

    foreach ($t in $Targets) {
        $r = $p.Range.Duplicate
        $found = $r.Find.Execute($t, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
        if (-not $found) {
            throw "Could not find '$t' in paragraph $idx (starting '$UniqueSubstring')"
        }
        $r.Font.Bold = 1
        $r.Font.Color = $highlightColor
    }
}

$plusMinus = [char]0x00B1

# "• Discovered systematic race coding errors ... from 23% to 64%"
Highlight-InParagraph "Discovered systematic race coding errors" @("23%", "64%")

# "• Utilized advanced sampling methods ... ±4.2% to ±2.1% ... 71% to 87% ..."
Highlight-InParagraph "Utilized advanced sampling methods to decrease survey margin of error from" @("${plusMinus}4.2%", "${plusMinus}2.1%", "71%", "87%")

# "• Trigonometric algorithm for boundary estimation ... 73.5% ... $4.7M ..."
Highlight-InParagraph "Trigonometric algorithm for boundary estimation" @("73.5%", "`$4.7M")

# "• Built real-time FEC analysis systems ... valued over $2 trillion"
Highlight-InParagraph "Built real-time FEC analysis systems" @("`$2")

# "• Predictive excellence: Utilized advanced sampling methods ... ±4.2% to ±2.1%"
Highlight-InParagraph "Predictive excellence" @("${plusMinus}4.2%", "${plusMinus}2.1%")

# "• Increased voter turnout prediction accuracy from 71% to 87%"
Highlight-InParagraph "Increased voter turnout prediction accuracy" @("71%", "87%")

# "• Methodological advancement: Improved segmentation accuracy 34% and survey incidence 28%"
Highlight-InParagraph "Methodological advancement" @("34%", "28%")
